$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A. This shifts the existing
# headers (Case ID, Trial Code, Arm, Arm Treatment, Diagnosis, Gender,
# Race, Ethnicity) one column to the right (B1:I1) and adds a new,
# still-empty, leading column A (the new "objects for trials" column).
$ws.Columns("A").Insert()

# The previous data row (the single sample trial record) is cleared out,
# leaving row 2 blank/empty - consistent with the new "objects" row
# having no data populated yet.
$ws.Range("B2:I2").ClearContents()

# Keep A1/A2 present as explicit (empty) cells, matching the shifted
# layout where column A has no header/value of its own yet.
$ws.Range("A1").Formula = '=""'
$ws.Range("A2").Formula = '=""'
